$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")

$ws.Range("A11").Value = "Country"
$ws.Range("B11").Value = "Country"
$ws.Range("C11").Value = "Land"
$ws.Range("D11").Value = "Need review"

# Replicate the explicit empty cell marker present in column E (e.g. <c r="E11"/>)
# by copying the formatting from D11 (default style) onto E11 without giving it a value.
$ws.Range("D11:E11").Copy()
$ws.Range("D11:E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
